$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.077.79'
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").Value = '1.651.96'
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5285'
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2598'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.46%  '
$ws.Range("E9").Value = '  +0.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.35'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07790'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.518'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.40%  '
$ws.Range("D13").Value = '1.658.65'
$ws.Range("E13").Value = '  -0.40%  '
$ws.Range("D14").Value = '1.879.21'
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5480'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8200'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.29'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.74%  '
$ws.Range("D18").Value = '26.077.60'
$ws.Range("E18").Value = '  -0.39%  '
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.584'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '190.83'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.38%  '
$ws.Range("E22").Value = '  +0.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.016'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.47%  '
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.98'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.02%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("E27").Value = '  -0.57%  '
$ws.Range("E28").Value = '  -1.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.452'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05786'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.46%  '
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.544'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.259'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.596'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.792'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.76%  '
$ws.Range("E36").Value = '  -0.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9445'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5743'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01613'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8546'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '104.15'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.001'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.708'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.96%  '
$ws.Range("D44").Value = '1.031.03'
$ws.Range("E44").Value = '  +2.76%  '
$ws.Range("D45").Value = '1.793.63'
$ws.Range("E45").Value = '  -0.48%  '
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("E47").Value = '  -0.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4326'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.862'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.41%  '
$ws.Range("E50").Value = '  -0.25%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.446'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.22%  '
